$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; B=1.459612070389937; C=0.04240448674262143; D=0.8054896365839992; E=8.660232485948974; G=10.96773867966553}
    @{Row=3; B=3.230985683306322; C=1.667794583268128; D=0.1575252929769615; E=0.496779210170732; G=5.553084769722144}
    @{Row=4; B=3.230985683306322; C=0.3127903958511391; D=0.1575252929769615; E=8.660232485948974; G=12.3615338580834}
    @{Row=5; B=0.003994804209775715; C=0.3127903958511391; D=3.900430680208489; E=8.660232485948974; G=12.87744836621838}
    @{Row=6; B=0.127881588408715; C=0.3127903958511391; D=0.1575252929769615; E=0.496779210170732; G=1.094976487407548}
    @{Row=7; B=3.230985683306322; C=1.667794583268128; D=337.1190423067083; E=8.660232485948974; G=350.6780550592317}
    @{Row=8; B=3.230985683306322; C=1.667794583268128; D=3.900430680208489; E=0.496779210170732; G=9.295990156953671}
    @{Row=9; B=0.6753301551942219; C=1.667794583268128; D=0.1575252929769615; E=8.660232485948974; G=11.16088251738829}
    @{Row=10; B=3.230985683306322; C=1.667794583268128; D=0.1575252929769615; E=0.496779210170732; G=5.553084769722144}
    @{Row=11; B=3.230985683306322; C=1.667794583268128; D=26.21740644021617; E=8.660232485948974; G=39.7764191927396}
    @{Row=12; B=3.230985683306322; C=10.29869402782916; D=26.21740644021617; E=8.660232485948974; G=48.40731863730063}
    @{Row=13; B=3.230985683306322; C=1.667794583268128; D=337.1190423067083; E=8.660232485948974; G=350.6780550592317}
    @{Row=14; B=1.459612070389937; C=1.667794583268128; D=9844.520545567508; E=8.660232485948974; G=9856.308184707115}
    @{Row=15; B=3.230985683306322; C=1.667794583268128; D=3.900430680208489; E=0.496779210170732; G=9.295990156953671}
    @{Row=16; B=3.230985683306322; C=1.667794583268128; D=0.1575252929769615; E=0.496779210170732; G=5.553084769722144}
    @{Row=17; B=0.127881588408715; C=1.667794583268128; D=3.900430680208489; E=8.660232485948974; G=14.35633933783431}
    @{Row=18; B=3.230985683306322; C=1.667794583268128; D=3.900430680208489; E=8.660232485948974; G=17.45944343273191}
    @{Row=19; B=3.230985683306322; C=1.667794583268128; D=0.8054896365839992; E=0.496779210170732; G=6.201049113329182}
    @{Row=20; B=1.459612070389937; C=1.667794583268128; D=0.8054896365839992; E=8.660232485948974; G=12.59312877619104}
    @{Row=21; B=0.127881588408715; C=0.3127903958511391; D=26.21740644021617; E=8.660232485948974; G=35.318310910425}
    @{Row=22; B=3.230985683306322; C=1.667794583268128; D=9844.520545567508; E=8.660232485948974; G=9858.079558320031}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 7).Value = $item.G
}
